$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '51.300.80'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.67%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.747.28'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.52%  '

# Row 4
$ws.Range("E4").Value = '  +0.07%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '115.13'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.29%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '332.74'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.00%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.531'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.51%  '

# Row 8
$ws.Range("E8").Value = '  +0.06%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.571'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +3.23%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '41.45'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.11%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.22'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.54%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0828'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.44%  '

# Row 13
$ws.Range("E13").Value = '  +2.77%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.64'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.73%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.177.28'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.68%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.748.21'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.81%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.889'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.54%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '51.217.99'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.68%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.82'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +4.84%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.01'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +4.02%  '

# Row 21
$ws.Range("E21").Value = '  +1.01%  '

# Row 22
$ws.Range("E22").Value = '  +0.30%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '279.93'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.59%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '70.19'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.28%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.64'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.53%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.92'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.05%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.999'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.04%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.37'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.43%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '36.08'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.45%  '

# Row 31
$ws.Range("E31").Value = '  -0.41%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '49.97'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.74%  '

# Row 33
$ws.Range("E33").Value = '  +2.53%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0829'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.22%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '19.62'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.36%  '

# Row 36
$ws.Range("B36").Value = 'ARBITRUM'
$ws.Range("C36").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.12'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.85%  '

# Row 37
$ws.Range("B37").Value = 'FirstDigitalUSD'
$ws.Range("C37").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.00'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.10%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.02'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.69%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.23'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.20%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '130.24'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +4.46%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '23.56'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.82%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0351'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +10.65%  '

# Row 43
$ws.Range("E43").Value = '  +4.21%  '

# Row 44
$ws.Range("E44").Value = '  +0.50%  '

# Row 45
$ws.Range("B45").Value = 'Stacks'
$ws.Range("C45").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.29'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +10.91%  '

# Row 46
$ws.Range("B46").Value = 'NEARProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.40'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.99%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.114.80'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.10%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.27'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.38%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.58'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.94%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.06'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.10%  '

# Row 51
$ws.Range("E51").Value = '  +9.30%  '
